# Translate outfit_detail_page.dart to English (and add main_screen.dart /
# outfit_detail_page.dart language rows), per commit:
#   "Dich outfit_detail_page.dart sang tieng Anh"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Fix the existing typo: "Favorite colours" -> "Favorite colors"
# ---------------------------------------------------------------------
$ws.Range("A45").Value2 = "Favorite colors"

# ---------------------------------------------------------------------
# 2) New section: main_screen.dart (rows 47-51)
# ---------------------------------------------------------------------

# Row 47 - section header (merged A47:B47), formatted like the other
# section headers (e.g. A32:B32).
[void]$ws.Range("A47:B47").Merge()
$ws.Rows.Item(47).RowHeight = 15.75
$ws.Range("A32:B32").Copy()
$ws.Range("A47:B47").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("A47").Value2 = "main_screen.dart"

# Row 48 - first data row, quote-prefixed style like A33:B33.
$ws.Rows.Item(48).RowHeight = 15.75
$ws.Range("A33:B33").Copy()
$ws.Range("A48:B48").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("A48").Value = "'Home"
$ws.Range("B48").Value = "'Trang chủ"

# Row 49 - regular style like A34:B34.
$ws.Rows.Item(49).RowHeight = 15.75
$ws.Range("A34:B34").Copy()
$ws.Range("A49:B49").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("A49").Value2 = "Closets"
$ws.Range("B49").Value2 = "Tủ đồ"

# Row 50 - alternate style like A35:B35.
$ws.Rows.Item(50).RowHeight = 15.75
$ws.Range("A35:B35").Copy()
$ws.Range("A50:B50").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("A50").Value2 = "Outfits"
$ws.Range("B50").Value2 = "Trang phục"

# Row 51 - alternate style like A35:B35.
$ws.Rows.Item(51).RowHeight = 15.75
$ws.Range("A35:B35").Copy()
$ws.Range("A51:B51").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("A51").Value2 = "Profile"
$ws.Range("B51").Value2 = "Cá nhân"

# ---------------------------------------------------------------------
# 3) New section: outfit_detail_page.dart (rows 53-55)
# ---------------------------------------------------------------------

# Row 53 - section header (merged A53:B53).
[void]$ws.Range("A53:B53").Merge()
$ws.Rows.Item(53).RowHeight = 15.75
$ws.Range("A32:B32").Copy()
$ws.Range("A53:B53").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("A53").Value2 = "outfit_detail_page.dart"

# Row 54 - first data row, quote-prefixed style like A33:B33.
$ws.Rows.Item(54).RowHeight = 15.75
$ws.Range("A33:B33").Copy()
$ws.Range("A54:B54").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("A54").Value = "'Fixed outfit"
$ws.Range("B54").Value = "'Trang phục cố định"

# Row 55 - regular style like A34:B34; taller row because of the longer
# wrapped text (matches row 16's 31.5pt pattern).
$ws.Rows.Item(55).RowHeight = 31.5
$ws.Range("A34:B34").Copy()
$ws.Range("A55:B55").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("A55").Value2 = "Items in this outfit are always worn together. Each item can only belong to one fixed outfit."
$ws.Range("B55").Value2 = "Các món đồ trong trang phục này luôn được mặc cùng nhau. Mỗi món đồ chỉ được thuộc một trang phục cố định"

# ---------------------------------------------------------------------
# 4) Restore the view state (scroll position / active selection) to
#    roughly match where a user editing this new content would land.
# ---------------------------------------------------------------------
[void]$ws.Activate()
[void]$ws.Range("A61").Select()
